$wb = $excel.ActiveWorkbook
$ws7 = $wb.Worksheets.Item("CO2 Capture")
$new = $wb.Worksheets.Add($null, $ws7)
$new.Name = "CLC Capture"
foreach ($s in $wb.Worksheets) {
    Write-Host $s.Index $s.Name
}
